$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Bulk updates to Price (D) and Volume(1h) (E) columns ---
$ws.Range("D2").Value = "35.386.22"
$ws.Range("E2").Value = "  +2.73%  "
$ws.Range("D3").Value = "1.842.34"
$ws.Range("E3").Value = "  +2.09%  "
$ws.Range("E4").Value = "  +0.27%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "229.95"
$ws.Range("E5").Value = "  +2.35%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.609"
$ws.Range("E6").Value = "  +3.14%  "
$ws.Range("E7").Value = "  +0.24%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "42.86"
$ws.Range("E8").Value = "  +12.10%  "
$ws.Range("E9").Value = "  +6.90%  "
$ws.Range("E10").Value = "  +3.30%  "
$ws.Range("E11").Value = "  +3.79%  "
$ws.Range("E12").Value = "  +2.14%  "
$ws.Range("D13").Value = "1.839.53"
$ws.Range("E13").Value = "  +1.93%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "11.28"
$ws.Range("E14").Value = "  +2.24%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.670"
$ws.Range("E15").Value = "  +7.00%  "
$ws.Range("E16").Value = "  +6.72%  "
$ws.Range("D17").Value = "35.367.61"
$ws.Range("E17").Value = "  +2.75%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "70.17"
$ws.Range("E18").Value = "  +3.39%  "
$ws.Range("D19").Value = "0.0₃0794"
$ws.Range("E19").Value = "  +3.38%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "244.50"
$ws.Range("E20").Value = "  +1.27%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.06"
$ws.Range("E21").Value = "  +9.03%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.64"
$ws.Range("E22").Value = "  +13.48%  "
$ws.Range("E23").Value = "  +0.31%  "
$ws.Range("E24").Value = "  +0.69%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "168.99"
$ws.Range("E25").Value = "  -1.13%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.90"
$ws.Range("E26").Value = "  +2.86%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.75"
$ws.Range("E27").Value = "  +2.06%  "
$ws.Range("E28").Value = "  +2.13%  "
$ws.Range("E29").Value = "  +13.59%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  +0.31%  "
$ws.Range("D31").Value = "3.256.31"
$ws.Range("E31").Value = "  +34.02%  "
$ws.Range("E32").Value = "  +6.35%  "
$ws.Range("E33").Value = "  +4.71%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.06"
$ws.Range("E34").Value = "  +5.87%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.86"
$ws.Range("E35").Value = "  +2.28%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "95.62"
$ws.Range("E36").Value = "  +16.05%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.685"
$ws.Range("E37").Value = "  +7.41%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.43"
$ws.Range("E40").Value = "  +5.83%  "
$ws.Range("E41").Value = "  +3.36%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.999"
$ws.Range("E42").Value = "  +6.12%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.26"
$ws.Range("E43").Value = "  +4.33%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.79"
$ws.Range("E44").Value = "  +8.10%  "
$ws.Range("E45").Value = "  +0.61%  "
$ws.Range("E46").Value = "  +0.05%  "
$ws.Range("E47").Value = "  +8.49%  "
$ws.Range("E48").Value = "  +1.82%  "
$ws.Range("D49").Value = "2.009.64"
$ws.Range("E49").Value = "  +2.22%  "
$ws.Range("E50").Value = "  +0.31%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "103.07"
$ws.Range("E51").Value = "  +0.94%  "


# --- Rows 38 and 39 swapped identity (Maker now ranks above TrustWalletToken) ---
$ws.Range("B38").Value = "Maker"
$ws.Range("C38").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D38").Value = "1.343.28"
$ws.Range("E38").Value = "  +1.62%  "

$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.08"
$ws.Range("E39").Value = "  +2.91%  "
